$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 (bold font, border, centered alignment) by
# copying H1's formatting into the new cells before setting their values.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J data columns for rows 2-30.
$ijValues = @{
    2  = @(1, 4)
    3  = @(6, 9)
    4  = @(2, 6)
    5  = @(1, 4)
    6  = @(6, 8)
    7  = @(6, 9)
    8  = @(2, 4)
    9  = @(7, 8)
    10 = @(7, 9)
    11 = @(6, 8)
    12 = @(2, 6)
    13 = @(1, 7)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 4)
    19 = @(1, 4)
    20 = @(3, 6)
    21 = @(6, 8)
    22 = @(9, 9)
    23 = @(7, 8)
    24 = @(6, 7)
    25 = @(6, 7)
    26 = @(6, 7)
    27 = @(1, 3)
    28 = @(1, 2)
    29 = @(4, 4)
    30 = @(3, 3)
}

foreach ($row in $ijValues.Keys) {
    $vals = $ijValues[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

Write-Output "applied I0/IF columns"
